# Estado de Cuenta - Comfenalco Cartagena
# "Elimina EC anteriores y se agregan nuevos, se modifica base de datos"
#
# Adds a new worker (YOLIMAR ELVIRA MARTINEZ MONTERO, CC 1044923399) to the
# account-statement table and updates the summary totals accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new row right below the existing worker (row 16), pushing the
#    footer/signature rows down by one, and copy row 16's formatting into it
#    so the new row keeps the same borders/number-formats as the table.
$ws.Rows("17:17").Insert()
$ws.Range("B16:J16").Copy()
$ws.Range("B17:J17").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# 2. Populate the new worker's data row.
$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "1044923399"
$ws.Range("D17").Value = "YOLIMAR ELVIRA MARTINEZ MONTERO"
$ws.Range("E17").Value = "2506"
$ws.Range("F17").Value = 46400
$ws.Range("G17").Value = 1160000

# 3. Update the summary block: total "Valor Mora" and worker count now
#    reflect both rows (52000 + 46400 = 98400; 2 trabajadores).
$ws.Range("E11").Value = 98400
$ws.Range("C13").Value = 2

# 4. Re-fit the data columns now that the table holds wider content
#    (longer name/id), matching Excel's automatic "best fit" column sizing.
$ws.Range("B1").ColumnWidth = 17.666666666666668
$ws.Range("C1").ColumnWidth = 15.833333333333334
$ws.Range("D1").ColumnWidth = 34.833333333333336
$ws.Range("E1").ColumnWidth = 12.666666666666666
$ws.Range("F1").ColumnWidth = 9.333333333333334
$ws.Range("G1").ColumnWidth = 13.5
$ws.Range("H1").ColumnWidth = 18.5
$ws.Range("I1").ColumnWidth = 17.333333333333332
$ws.Range("J1").ColumnWidth = 14.166666666666666
